$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.133.97'
$ws.Range("E2").Value = '  -0.96%  '

$ws.Range("D3").Value = '1.670.23'
$ws.Range("E3").Value = '  -1.39%  '

$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = '  -0.72%  '

$ws.Range("D5").Value = "'210.49"
$ws.Range("E5").Value = '  -4.01%  '

$ws.Range("D6").Value = "'0.5243"
$ws.Range("E6").Value = '  -4.69%  '

$ws.Range("E7").Value = '  -0.71%  '

$ws.Range("D8").Value = "'0.2649"
$ws.Range("E8").Value = '  -3.33%  '

$ws.Range("D9").Value = "'0.06276"
$ws.Range("E9").Value = '  -2.92%  '

$ws.Range("D10").Value = "'21.14"
$ws.Range("E10").Value = '  -3.94%  '

$ws.Range("D11").Value = "'0.07513"
$ws.Range("E11").Value = '  -1.98%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.662.81'
$ws.Range("E12").Value = '  -2.12%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'4.439"
$ws.Range("E13").Value = '  -2.23%  '

$ws.Range("D14").Value = "'0.5621"
$ws.Range("E14").Value = '  -3.63%  '

$ws.Range("D15").Value = "'0.000008017"
$ws.Range("E15").Value = '  -4.12%  '

$ws.Range("D16").Value = "'66.35"
$ws.Range("E16").Value = '  +1.30%  '

$ws.Range("D17").Value = '26.188.08'
$ws.Range("E17").Value = '  -0.90%  '

$ws.Range("E18").Value = '  -0.77%  '

$ws.Range("D19").Value = "'4.796"
$ws.Range("E19").Value = '  -2.86%  '

$ws.Range("D20").Value = "'187.55"
$ws.Range("E20").Value = '  -2.27%  '

$ws.Range("D21").Value = "'10.36"
$ws.Range("E21").Value = '  -5.56%  '

$ws.Range("D22").Value = "'6.175"
$ws.Range("E22").Value = '  -1.19%  '

$ws.Range("E23").Value = '  -0.68%  '

$ws.Range("D24").Value = "'148.19"
$ws.Range("E24").Value = '  -0.45%  '

$ws.Range("D25").Value = "'0.1248"
$ws.Range("E25").Value = '  -5.73%  '

$ws.Range("D26").Value = "'7.591"
$ws.Range("E26").Value = '  -4.06%  '

$ws.Range("D27").Value = "'15.95"
$ws.Range("E27").Value = '  +1.19%  '

$ws.Range("D28").Value = "'0.06216"
$ws.Range("E28").Value = '  -0.99%  '

$ws.Range("E29").Value = '  -1.68%  '

$ws.Range("E30").Value = '  -4.06%  '

$ws.Range("D31").Value = "'3.469"

$ws.Range("D32").Value = "'3.432"
$ws.Range("E32").Value = '  -4.64%  '

$ws.Range("D33").Value = "'1.622"
$ws.Range("E33").Value = '  -3.68%  '

$ws.Range("D34").Value = "'0.9942"
$ws.Range("E34").Value = '  -4.62%  '

$ws.Range("D35").Value = "'0.6037"
$ws.Range("E35").Value = '  -1.81%  '

$ws.Range("D36").Value = "'2.402"
$ws.Range("E36").Value = '  -0.44%  '

$ws.Range("D37").Value = "'2.714"
$ws.Range("E37").Value = '  +0.14%  '

$ws.Range("D38").Value = "'6.111"
$ws.Range("E38").Value = '  -1.45%  '

$ws.Range("D39").Value = "'0.01612"
$ws.Range("E39").Value = '  -1.87%  '

$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '1.073.54'
$ws.Range("E40").Value = '  -3.85%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = "'0.8668"
$ws.Range("E41").Value = '  -2.42%  '

$ws.Range("D42").Value = "'1.005"
$ws.Range("E42").Value = '  -1.11%  '

$ws.Range("D43").Value = "'99.99"
$ws.Range("E43").Value = '  -1.85%  '

$ws.Range("E45").Value = '  +1.20%  '

$ws.Range("D46").Value = "'56.02"
$ws.Range("E46").Value = '  -2.55%  '

$ws.Range("E47").Value = '  -1.56%  '

$ws.Range("D48").Value = "'0.05237"
$ws.Range("E48").Value = '  -0.98%  '

$ws.Range("D49").Value = "'7.964"
$ws.Range("E49").Value = '  -2.65%  '

$ws.Range("D50").Value = "'0.4255"
$ws.Range("E50").Value = '  -1.13%  '

$ws.Range("D51").Value = "'5.990"
$ws.Range("E51").Value = '  -1.95%  '
